# excel for multiple product upload done
#
# Renames the two sheets, restructures the product-catalog header row to the
# new camelCase schema (with three/two repeated color+size:quantity+photo
# "variant" groups), and rewrites the sample product rows to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename sheets
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("First")
$ws1.Name = "FirstSheet"

$ws2 = $wb.Worksheets.Item("Second")
$ws2.Name = "SecondSheet"

# ---------------------------------------------------------------------------
# 2) FirstSheet — header row + 3 product rows, 22 columns (A:V)
# ---------------------------------------------------------------------------
$header1 = @("name", "code", "subTitle", "description", "price", "costPrice", "category", "gender", "brand", "weight", "composition", "tags", "orderQuantity", "color", "size:quantity", "photo", "color", "size:quantity", "photo", "color", "size:quantity", "photo")

$row2 = @("Product 1", "XYZ", "Subtitle A", "This is product 1", 1199.0, 999.0, "Kurti", "Female", "Brand X", "8 kg", "100% Cotton", "Tag1, Tag3", "7, 8, 1", "#FF0000", "X:23, L:40", "link1, link2, link3", "#f5cb42", "XS:100, M:50, L:80", "link1, link2, link3", $null, $null, $null)
$row3 = @("Product 2", "ABC", "Subtitle B", "This is product 2", 560.0, 450.0, "Shirt", "Female", "Brand Y", "1 kg", "100% Cotton", "Tag2, Tag4", "1, 4, 9", "#00FF00", "XS:60, M:50, L:20", "link1, link2", $null, $null, $null, $null, $null, $null)
$row4 = @("Product 3", "AYD", "Subtitle B", "This is product 3", 1599.0, 1300.0, "T-shirt", "Male", "Brand X", "7 kg", "100% Cotton", "Tag2, Tag3", "5, 6, 2", "#00FFFF", "L:100, XL:45", "link1, link2, link3, link4", "#a83275", "L:75, XS:60, XL:25", "link1, link2, link3, link4", "#c92926", "M:40, L:60, XXL:70", "link1, link2, link3, link4")

for ($i = 0; $i -lt $header1.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $header1[$i]
}
for ($i = 0; $i -lt $row2.Length; $i++) {
    if ($null -ne $row2[$i]) { $ws1.Cells.Item(2, $i + 1).Value = $row2[$i] }
}
for ($i = 0; $i -lt $row3.Length; $i++) {
    if ($null -ne $row3[$i]) { $ws1.Cells.Item(3, $i + 1).Value = $row3[$i] }
}
for ($i = 0; $i -lt $row4.Length; $i++) {
    if ($null -ne $row4[$i]) { $ws1.Cells.Item(4, $i + 1).Value = $row4[$i] }
}

# ---------------------------------------------------------------------------
# 3) SecondSheet — header row + a single product row, 17 columns (A:Q);
#    rows 3 and 4 (which used to hold Product 2 / Product 3 duplicates) are
#    cleared out since the sheet now only carries one sample row.
# ---------------------------------------------------------------------------
$header2 = @("name", "code", "subTitle", "description", "category", "gender", "brand", "weight", "composition", "tags", "orderQuantity", "color", "size:quantity", "photo", "color", "size:quantity", "photo")
$s2row2 = @("Product 4", "OKK", "Subtitle K", "This is product 4", "Jeans", "Male", "Brand OK", "8 kg", "100% Cotton", "Tag1, Tag3", "7, 8, 1", "#FF0000", "X:23, L:40", "link1, link2, link3", "#f5cb42", "XS:100, M:50, L:80", "link1, link2, link3")

for ($i = 0; $i -lt $header2.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $header2[$i]
}
for ($i = 0; $i -lt $s2row2.Length; $i++) {
    if ($null -ne $s2row2[$i]) { $ws2.Cells.Item(2, $i + 1).Value = $s2row2[$i] }
}

$ws2.Rows.Item(3).Clear()
$ws2.Rows.Item(4).Clear()
$ws2.Rows.Item(3).RowHeight = 15.75
$ws2.Rows.Item(4).RowHeight = 15.75

Write-Output "done"
